{"js": "const replacements = [\n  [\"2025-12-27 Saturday\", \"2025-12-28 Sunday\"],\n  [\"489\u00d74=\", \"905\u00d73=\"],\n  [\"434\u00d77=\", \"883\u00d72=\"],\n  [\"212\u00d77=\", \"426\u00d76=\"],\n  [\"780\u00d76=\", \"365\u00d73=\"],\n  [\"881\u00d77=\", \"412\u00d79=\"],\n  [\"650\u00d77=\", \"913\u00d79=\"],\n  [\"376\u00d77=\", \"290\u00d77=\"],\n  [\"925\u00d73=\", \"907\u00d77=\"],\n  [\"755\u00d73=\", \"699\u00d79=\"],\n  [\"806\u00d76=\", \"863\u00d78=\"],\n  [\"974\u00d76=\", \"173\u00d79=\"],\n  [\"626\u00d76=\", \"269\u00d79=\"],\n  [\"261\u00d76=\", \"592\u00d74=\"],\n  [\"457\u00d77=\", \"263\u00d74=\"],\n  [\"541\u00d76=\", \"750\u00d72=\"],\n  [\"975\u00d78=\", \"671\u00d79=\"],\n  [\"961\u00d73=\", \"421\u00d76=\"],\n  [\"796\u00d72=\", \"712\u00d75=\"],\n  [\"914\u00d74=\", \"370\u00d78=\"],\n  [\"467\u00d77=\", \"995\u00d76=\"],\n  [\"136\u00d73=\", \"512\u00d73=\"],\n  [\"693\u00d79=\", \"330\u00d73=\"],\n  [\"648\u00d75=\", \"847\u00d78=\"],\n  [\"406\u00d76=\", \"377\u00d75=\"],\n  [\"829\u00d73=\", \"170\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"2025-12-27 Saturday\"; new = \"2025-12-28 Sunday\"},\n    @{old = \"489\u00d74=\"; new = \"905\u00d73=\"},\n    @{old = \"434\u00d77=\"; new = \"883\u00d72=\"},\n    @{old = \"212\u00d77=\"; new = \"426\u00d76=\"},\n    @{old = \"780\u00d76=\"; new = \"365\u00d73=\"},\n    @{old = \"881\u00d77=\"; new = \"412\u00d79=\"},\n    @{old = \"650\u00d77=\"; new = \"913\u00d79=\"},\n    @{old = \"376\u00d77=\"; new = \"290\u00d77=\"},\n    @{old = \"925\u00d73=\"; new = \"907\u00d77=\"},\n    @{old = \"755\u00d73=\"; new = \"699\u00d79=\"},\n    @{old = \"806\u00d76=\"; new = \"863\u00d78=\"},\n    @{old = \"974\u00d76=\"; new = \"173\u00d79=\"},\n    @{old = \"626\u00d76=\"; new = \"269\u00d79=\"},\n    @{old = \"261\u00d76=\"; new = \"592\u00d74=\"},\n    @{old = \"457\u00d77=\"; new = \"263\u00d74=\"},\n    @{old = \"541\u00d76=\"; new = \"750\u00d72=\"},\n    @{old = \"975\u00d78=\"; new = \"671\u00d79=\"},\n    @{old = \"961\u00d73=\"; new = \"421\u00d76=\"},\n    @{old = \"796\u00d72=\"; new = \"712\u00d75=\"},\n    @{old = \"914\u00d74=\"; new = \"370\u00d78=\"},\n    @{old = \"467\u00d77=\"; new = \"995\u00d76=\"},\n    @{old = \"136\u00d73=\"; new = \"512\u00d73=\"},\n    @{old = \"693\u00d79=\"; new = \"330\u00d73=\"},\n    @{old = \"648\u00d75=\"; new = \"847\u00d78=\"},\n    @{old = \"406\u00d76=\"; new = \"377\u00d75=\"},\n    @{old = \"829\u00d73=\"; new = \"170\u00d79=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
